# Remove the "HFE" (Huainan, China) colo row from the data center colocation
# list. This deletes the entire worksheet row 261, shifting every
# subsequent row (XFN, XNN, DAD, JXG, ... YHZ) up by one and shrinking the
# used range from A1:H334 to A1:H333.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(261).Delete()
